# The sheet currently has headers in B1:T1:
#   B1=fullName, C1=level, D1=teamName, E1=Loại, F1=KR phòng, G1=KR team,
#   H1=KR cá nhân, I1=Công thức tính, J1=Nguồn dữ liệu, K1=Định kỳ tính,
#   L1=Đơn vị tính, M1=Điều kiện, N1=Norm, O1=% Trọng số chỉ tiêu,
#   P1=Kết quả, Q1=Tỷ lệ, R1=Tổng thời gian dự kiến/ ước tính công việc (giờ),
#   S1=Tổng thời gian thực hiện công việc thực tế (giờ), T1=Note
#
# Target: insert a new "Name" column right after fullName (which becomes
# employeeId), and a new "krId" column right after "Loại" (before "KR phòng").
# Everything else shifts right by these two new columns, ending at V1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new columns. Insert the earlier (left-most) one first, so
# that the column letter used for the second insert already accounts for
# the shift caused by the first insert.
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("G1").EntireColumn.Insert()

# Rename the existing fullName header to employeeId, and fill the new cells.
$ws.Range("B1").Value = "employeeId"
$ws.Range("C1").Value = "Name"
$ws.Range("G1").Value = "krId"

# Match header formatting (bold, centered, bordered) used by the rest of row 1.
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("C1").Borders.LineStyle = 1

$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").VerticalAlignment = -4160
$ws.Range("G1").Borders.LineStyle = 1
